# updated Day 6 training materials
#
# 1) The "datetimeFigureOut" date fields on the slide master, every slide
#    layout and the notes master get their cached text bumped from
#    8/6/2023 to 8/10/2023 (PowerPoint re-caches these automatically when
#    the deck is opened/saved on a later day).
# 2) Slide 2's big "Assisted Practice" caption becomes "Demo", rendered in
#    Microsoft Sans Serif (instead of Verdana) with tightened character
#    spacing.

$p = $ppt.ActivePresentation

function Find-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
        }
        if ($isDate) { return $sh }
    }
    return $null
}

function Update-DateShape($sh) {
    if ($sh -eq $null) { return }
    $tr = $sh.TextFrame.TextRange
    if ($tr.Text -eq "8/6/2023") {
        $tr.Text = "8/10/2023"
    }
}

# --- Slide master ---
$master = $p.SlideMaster
Update-DateShape (Find-DatePlaceholder $master.Shapes)

# --- Every slide layout hanging off the master ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape (Find-DatePlaceholder $layout.Shapes)
}

# --- Notes master ---
Update-DateShape (Find-DatePlaceholder $p.NotesMaster.Shapes)

# --- Slide 2: "Assisted Practice" -> "Demo" ---
$slide = $p.Slides.Item(2)
$shape = $slide.Shapes.Item(4)
$heightPts = $shape.Height

$textRange = $shape.TextFrame.TextRange
$textRange.Text = "Demo"
$textRange.Font.Name = "Microsoft Sans Serif"
$textRange.Font.NameComplexScript = "Microsoft Sans Serif"
$textRange.Font.Spacing = -0.8

# Changing the run's font nudges the autosized textbox height; put it back.
$shape.Height = $heightPts
